# Auto-generated edit script: apply scheduled-runner market-price updates
# to the Leve profit-calculation columns (H:N) across all profession sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1317.0625
$ws.Range("I19").Value = 1669.1111
$ws.Range("J19").Value = 864.4286
$ws.Range("K19").Value = 1669.1111
$ws.Range("L19").Value = 864.4286
$ws.Range("M19").Value = -1494.1111
$ws.Range("N19").Value = -1214.4286
$ws.Range("H40").Value = 2926509.5
$ws.Range("J40").Value = 10104661
$ws.Range("L40").Value = 10104661
$ws.Range("N40").Value = -10105011
$ws.Range("H70").Value = 6488.0586
$ws.Range("I70").Value = 2527.182
$ws.Range("J70").Value = 13749.667
$ws.Range("K70").Value = 7581.545999999999
$ws.Range("L70").Value = 41249.001
$ws.Range("M70").Value = -7311.545999999999
$ws.Range("N70").Value = -41789.001
$ws.Range("H73").Value = 6488.0586
$ws.Range("I73").Value = 2527.182
$ws.Range("J73").Value = 13749.667
$ws.Range("K73").Value = 7581.545999999999
$ws.Range("L73").Value = 41249.001
$ws.Range("M73").Value = -6645.545999999999
$ws.Range("N73").Value = -43121.001
$ws.Range("H98").Value = 2759.1667
$ws.Range("I98").Value = 2759.1667
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 2759.1667
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -1261.1667
$ws.Range("N98").ClearContents()
$ws.Range("H122").Value = 2759.1667
$ws.Range("I122").Value = 2759.1667
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8277.500100000001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5827.500100000001
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 89999
$ws.Range("J7").Value = 89999
$ws.Range("L7").Value = 89999
$ws.Range("N7").Value = -90227
$ws.Range("H45").Value = 1811.75
$ws.Range("I45").Value = 1649.3334
$ws.Range("K45").Value = 1649.3334
$ws.Range("M45").Value = -1272.3334
$ws.Range("H122").Value = 4611.4443
$ws.Range("I122").Value = 2786.2856
$ws.Range("K122").Value = 8358.856800000001
$ws.Range("M122").Value = -5908.856800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 688.1667
$ws.Range("I19").Value = 532.25
$ws.Range("J19").Value = 1000
$ws.Range("K19").Value = 532.25
$ws.Range("L19").Value = 1000
$ws.Range("M19").Value = -362.25
$ws.Range("N19").Value = -1340
$ws.Range("H24").Value = 688.1667
$ws.Range("I24").Value = 532.25
$ws.Range("J24").Value = 1000
$ws.Range("K24").Value = 532.25
$ws.Range("L24").Value = 1000
$ws.Range("M24").Value = -362.25
$ws.Range("N24").Value = -1340
$ws.Range("H95").Value = 27999
$ws.Range("J95").Value = 27999
$ws.Range("L95").Value = 27999
$ws.Range("N95").Value = -33491
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H97").Value = 39998.5
$ws.Range("J97").Value = 39998.5
$ws.Range("L97").Value = 39998.5
$ws.Range("N97").Value = -41980.5
$ws.Range("H103").Value = 33498
$ws.Range("I103").Value = 26197.8
$ws.Range("K103").Value = 26197.8
$ws.Range("M103").Value = -25025.8
$ws.Range("H104").Value = 17333.334
$ws.Range("I104").Value = 12000
$ws.Range("J104").Value = 20000
$ws.Range("K104").Value = 12000
$ws.Range("L104").Value = 20000
$ws.Range("M104").Value = -9379
$ws.Range("N104").Value = -25242
$ws.Range("H141").Value = 302813.62
$ws.Range("I141").Value = 100000
$ws.Range("J141").Value = 331787
$ws.Range("K141").Value = 100000
$ws.Range("L141").Value = 331787
$ws.Range("M141").Value = -94820
$ws.Range("N141").Value = -342147

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 81.545456
$ws.Range("I2").Value = 40.625
$ws.Range("J2").Value = 104.92857
$ws.Range("K2").Value = 243.75
$ws.Range("L2").Value = 629.57142
$ws.Range("M2").Value = -130.75
$ws.Range("N2").Value = -855.57142
$ws.Range("H23").Value = 1994.5
$ws.Range("I23").Value = 1994.5
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 5983.5
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -5748.5
$ws.Range("N23").ClearContents()
$ws.Range("H37").Value = 137568
$ws.Range("J37").Value = 137568
$ws.Range("L37").Value = 412704
$ws.Range("N37").Value = -412928

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()
$ws.Range("H63").Value = 40103
$ws.Range("I63").Value = 40103
$ws.Range("K63").Value = 40103
$ws.Range("M63").Value = -39417
$ws.Range("H66").Value = 40103
$ws.Range("I66").Value = 40103
$ws.Range("K66").Value = 120309
$ws.Range("M66").Value = -116877
$ws.Range("H80").Value = 1915.75
$ws.Range("J80").Value = 2433.3333
$ws.Range("L80").Value = 2433.3333
$ws.Range("N80").Value = -4429.3333
$ws.Range("H83").Value = 1915.75
$ws.Range("J83").Value = 2433.3333
$ws.Range("L83").Value = 12166.6665
$ws.Range("N83").Value = -22150.6665
$ws.Range("H122").Value = 3268
$ws.Range("I122").Value = 1695.7
$ws.Range("K122").Value = 5087.1
$ws.Range("M122").Value = -2637.1
$ws.Range("H132").Value = 8336644.5
$ws.Range("I132").Value = 8932063
$ws.Range("J132").Value = 789
$ws.Range("K132").Value = 26796189
$ws.Range("L132").Value = 2367
$ws.Range("M132").Value = -26793659
$ws.Range("N132").Value = -7427

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 10250
$ws.Range("I22").Value = 10500
$ws.Range("K22").Value = 10500
$ws.Range("M22").Value = -10205
$ws.Range("H27").Value = 10250
$ws.Range("I27").Value = 10500
$ws.Range("K27").Value = 10500
$ws.Range("M27").Value = -10393
$ws.Range("H40").Value = 2822.3333
$ws.Range("I40").Value = 2806.8
$ws.Range("K40").Value = 2806.8
$ws.Range("M40").Value = -2670.8
$ws.Range("H68").Value = 8334499.5
$ws.Range("I68").Value = 8334499.5
$ws.Range("K68").Value = 8334499.5
$ws.Range("M68").Value = -8333750.5
$ws.Range("H71").Value = 8334499.5
$ws.Range("I71").Value = 8334499.5
$ws.Range("K71").Value = 41672497.5
$ws.Range("M71").Value = -41668753.5
$ws.Range("H82").Value = 709.6
$ws.Range("J82").Value = 374
$ws.Range("L82").Value = 374
$ws.Range("N82").Value = -1096
$ws.Range("H85").Value = 709.6
$ws.Range("J85").Value = 374
$ws.Range("L85").Value = 374
$ws.Range("N85").Value = -2870
$ws.Range("H99").Value = 29424.75
$ws.Range("I99").Value = 25900
$ws.Range("K99").Value = 25900
$ws.Range("M99").Value = -22905
$ws.Range("H102").Value = 69780
$ws.Range("J102").Value = 69780
$ws.Range("L102").Value = 69780
$ws.Range("N102").Value = -76270
$ws.Range("H122").Value = 7682.1665
$ws.Range("I122").Value = 8017.1816
$ws.Range("J122").Value = 3997
$ws.Range("K122").Value = 24051.5448
$ws.Range("L122").Value = 11991
$ws.Range("M122").Value = -21601.5448
$ws.Range("N122").Value = -16891
$ws.Range("H129").Value = 60000
$ws.Range("J129").Value = 60000
$ws.Range("L129").Value = 60000
$ws.Range("N129").Value = -70000

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 49325.5
$ws.Range("J27").Value = 49325.5
$ws.Range("L27").Value = 49325.5
$ws.Range("N27").Value = -49463.5
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H102").Value = 66337
$ws.Range("J102").Value = 66337
$ws.Range("L102").Value = 66337
$ws.Range("N102").Value = -72827
$ws.Range("H106").Value = 81438.25
$ws.Range("J106").Value = 81438.25
$ws.Range("L106").Value = 81438.25
$ws.Range("N106").Value = -83962.25
$ws.Range("H115").Value = 49376
$ws.Range("J115").Value = 49376
$ws.Range("L115").Value = 49376
$ws.Range("N115").Value = -52510
$ws.Range("H132").Value = 12199789
$ws.Range("I132").Value = 19232974
$ws.Range("K132").Value = 57698922
$ws.Range("M132").Value = -57696392
